# "Tried to implement Penality Reward System (unfinished)"
#
# This script reproduces the target edit across the workbook's five sheets:
#
#  - "Weekly Sales": append one more weekly bucket (row 63) with y = 0.
#  - "Daily PO": clear out all PO rows (the underlying PO data was wiped)
#       and drop the now-unused "Index" column (AG).
#  - "Merged (Optional)": the first two (now-cancelled) PO-linked rows are
#       dropped, every later row shifts up accordingly, and the series is
#       extended by one more week (row 63) to stay in sync with
#       "Weekly Sales".
#  - "PO Volume Insights" / "PO Prediction": zero out the PO volume
#       aggregates / forecast, since there is no PO data left to derive
#       them from.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Weekly Sales" -- append a new weekly row.
# ---------------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Sales")

$wsWeekly.Range("A63").Value = 45662.99999999999
$wsWeekly.Range("A63").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsWeekly.Range("B63").Value = 0

# ---------------------------------------------------------------------------
# 2) "Daily PO" -- remove all PO data rows and the "Index" column (AG).
# ---------------------------------------------------------------------------
$wsDailyPO = $wb.Worksheets.Item("Daily PO")

# Rows 2 and 3 hold the only two PO records; delete both (deleting the same
# row index twice removes both original rows 2 and 3).
$wsDailyPO.Rows.Item(2).Delete()
$wsDailyPO.Rows.Item(2).Delete()

# Column AG ("Index") is no longer used once the PO rows are gone.
$wsDailyPO.Columns.Item(33).Delete()

# ---------------------------------------------------------------------------
# 3) "Merged (Optional)" -- drop the first two PO-era rows (old rows 7 & 8)
#    and extend the series with a new trailing week (row 63).
# ---------------------------------------------------------------------------
$wsMerged = $wb.Worksheets.Item("Merged (Optional)")

$wsMerged.Rows.Item(7).Delete()
$wsMerged.Rows.Item(7).Delete()

$wsMerged.Range("A63").Value = 45662.99999999999
$wsMerged.Range("A63").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsMerged.Range("B63").Value = 0
$wsMerged.Range("C63").Value = 0

# ---------------------------------------------------------------------------
# 4) "PO Volume Insights" -- zero out the aggregate PO volume stats.
# ---------------------------------------------------------------------------
$wsInsights = $wb.Worksheets.Item("PO Volume Insights")

$wsInsights.Range("A2").Value = 0
$wsInsights.Range("B2").Value = 0
$wsInsights.Range("C2").Value = 0
$wsInsights.Range("D2").Value = 0

# ---------------------------------------------------------------------------
# 5) "PO Prediction" -- zero out the predicted next daily PO quantity.
# ---------------------------------------------------------------------------
$wsPrediction = $wb.Worksheets.Item("PO Prediction")

$wsPrediction.Range("A2").Value = 0
